$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14 - new shared string "Thiếu thêm, xóa, sửa" created first (index 33)
$ws.Range("H14").Value = 0.25
$ws.Range("J14").Value = "Thiếu thêm, xóa, sửa"

# Row 15 - reuses the same shared string
$ws.Range("H15").Value = 0.25
$ws.Range("J15").Value = "Thiếu thêm, xóa, sửa"

# Row 12 - new shared string "Chưa edit + update ảnh" created next (index 34)
$ws.Range("H12").Value = 0.5
$ws.Range("J12").Value = "Chưa edit + update ảnh"

# Row 13
$ws.Range("H13").Value = 1

# Update selection to match new active cell
$ws.Range("H13").Select()
